$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248 }
    3 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043;  E = 13.86384647080068;   G = 18.91276827552123 }
    4 = @{ B = 0.2881169905109251; C = 0.04103571897497393; D = 0.1496068669990043; E = 0.5333859586016987;  G = 1.012145535086602 }
    5 = @{ B = 0.6545652718822623; C = 2919.202174992006;   D = 0.7210945179870265; E = 13.86384647080068;   G = 2934.441681252676 }
    6 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987;  G = 4.327115817150455 }
    7 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987;  G = 8.656069925401464 }
    8 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265;  E = 0.5333859586016987;  G = 4.327115817150455 }
    9 = @{ B = 0.1169995834814548; C = 0.3048912486333797; D = 0.7210945179870265;  E = 0.5333859586016987;  G = 1.67637130870356 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
